$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-SceneRow($r, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j) {
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
    $ws.Range("I$r").Value = $i
    $ws.Range("J$r").Value = $j
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("J$r").NumberFormat = "@"
}

# New rows for camera offset position / rotation configuration
Set-SceneRow 13 "CamOffestPos" "string" $false $false $false $true 0 0 "Friend" "acctorid"
Set-SceneRow 14 "CamOffestRot" "string" $false $false $false $true 0 0 "Friend" "acctorid"

# Keep the TRUE/FALSE list validation applied to the new F column cells.
# Re-create the rule on F13:F14 explicitly so the sqref of the existing
# rule is split the same way Excel splits it when editing cells that
# already fall inside a previously-saved validated range.
$ws.Range("F13:F14").Validation.Delete()
$ws.Range("F13:F14").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Match the selection Excel leaves behind after editing the last row
$ws.Range("A14").Select()
